$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4511.9375
$ws.Range("I40").Value = 5780.6665
$ws.Range("K40").Value = 5780.6665
$ws.Range("M40").Value = -5605.6665
# Row 74
$ws.Range("H74").Value = 14459.348
$ws.Range("I74").Value = 14450.842
$ws.Range("K74").Value = 14450.842
$ws.Range("M74").Value = -13514.842
# Row 77
$ws.Range("H77").Value = 14459.348
$ws.Range("I77").Value = 14450.842
$ws.Range("K77").Value = 72254.21000000001
$ws.Range("M77").Value = -67574.21000000001
# Row 100
$ws.Range("H100").Value = 3895
$ws.Range("I100").Value = 3895
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3895
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3354
$ws.Range("N100").ClearContents()
# Row 121
$ws.Range("H121").Value = 2331.6667
$ws.Range("J121").Value = 2331.6667
$ws.Range("L121").Value = 6995.000100000001
$ws.Range("N121").Value = -10489.0001
# Row 138
$ws.Range("H138").Value = 265996.9
$ws.Range("I138").Value = 3818.8333
$ws.Range("J138").Value = 354620.47
$ws.Range("K138").Value = 11456.4999
$ws.Range("L138").Value = 1063861.41
$ws.Range("M138").Value = -6316.499899999999
$ws.Range("N138").Value = -1074141.41

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1052.619
$ws.Range("J2").Value = 1187.25
$ws.Range("L2").Value = 1187.25
$ws.Range("N2").Value = -1413.25
# Row 45
$ws.Range("H45").Value = 65231.43
$ws.Range("I45").Value = 65231.43
$ws.Range("K45").Value = 65231.43
$ws.Range("M45").Value = -64854.43
# Row 61
$ws.Range("H61").Value = 5246.1
$ws.Range("I61").Value = 3089.2
$ws.Range("J61").Value = 7403
$ws.Range("K61").Value = 3089.2
$ws.Range("L61").Value = 7403
$ws.Range("M61").Value = -2877.2
$ws.Range("N61").Value = -7827
# Row 74
$ws.Range("H74").Value = 193584.34
$ws.Range("I74").Value = 278881.8
$ws.Range("K74").Value = 278881.8
$ws.Range("M74").Value = -278007.8
# Row 77
$ws.Range("H77").Value = 193584.34
$ws.Range("I77").Value = 278881.8
$ws.Range("K77").Value = 1394409
$ws.Range("M77").Value = -1390041
# Row 110
$ws.Range("H110").Value = 3718.3333
$ws.Range("I110").Value = 1854
$ws.Range("K110").Value = 1854
$ws.Range("M110").Value = 191
# Row 116
$ws.Range("H116").Value = 1052.619
$ws.Range("J116").Value = 1187.25
$ws.Range("L116").Value = 1187.25
$ws.Range("N116").Value = -5775.25
# Row 124
$ws.Range("H124").Value = 54666.332
$ws.Range("J124").Value = 54666.332
$ws.Range("L124").Value = 54666.332
$ws.Range("N124").Value = -64486.332
# Row 132
$ws.Range("H132").Value = 4248.1904
$ws.Range("I132").Value = 1979.5714
$ws.Range("K132").Value = 5938.7142
$ws.Range("M132").Value = -3408.7142
# Row 136
$ws.Range("H136").Value = 5246.1
$ws.Range("I136").Value = 3089.2
$ws.Range("J136").Value = 7403
$ws.Range("K136").Value = 9267.599999999999
$ws.Range("L136").Value = 22209
$ws.Range("M136").Value = -6717.599999999999
$ws.Range("N136").Value = -27309

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1052.619
$ws.Range("J3").Value = 1187.25
$ws.Range("L3").Value = 1187.25
$ws.Range("N3").Value = -1415.25
# Row 6
$ws.Range("H6").Value = 10696.5
$ws.Range("J6").Value = 10696.5
$ws.Range("L6").Value = 10696.5
$ws.Range("N6").Value = -10922.5
# Row 86
$ws.Range("H86").Value = 3247.7273
$ws.Range("I86").Value = 2922.5
$ws.Range("K86").Value = 2922.5
$ws.Range("M86").Value = -1799.5
# Row 89
$ws.Range("H89").Value = 3247.7273
$ws.Range("I89").Value = 2922.5
$ws.Range("K89").Value = 14612.5
$ws.Range("M89").Value = -8996.5
# Row 99
$ws.Range("H99").Value = 94305
$ws.Range("I99").Value = 102985.5
$ws.Range("K99").Value = 102985.5
$ws.Range("M99").Value = -101487.5
# Row 134
$ws.Range("H134").Value = 7999.5
$ws.Range("I134").Value = 7998
$ws.Range("K134").Value = 23994
$ws.Range("M134").Value = -21459

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1710.84
$ws.Range("J16").Value = 1803.5
$ws.Range("L16").Value = 1803.5
$ws.Range("N16").Value = -2377.5
# Row 58
$ws.Range("H58").Value = 3165.4736
$ws.Range("I58").Value = 1946.4
$ws.Range("K58").Value = 1946.4
$ws.Range("M58").Value = -1743.4
# Row 62
$ws.Range("H62").Value = 16681832
$ws.Range("J62").Value = 22664
$ws.Range("L62").Value = 22664
$ws.Range("N62").Value = -23912
# Row 65
$ws.Range("H65").Value = 16681832
$ws.Range("J65").Value = 22664
$ws.Range("L65").Value = 113320
$ws.Range("N65").Value = -119560
# Row 113
$ws.Range("H113").Value = 1710.84
$ws.Range("J113").Value = 1803.5
$ws.Range("L113").Value = 1803.5
$ws.Range("N113").Value = -6143.5
# Row 134
$ws.Range("H134").Value = 3164.7878
$ws.Range("I134").Value = 2942.4075
$ws.Range("K134").Value = 8827.2225
$ws.Range("M134").Value = -6292.2225
# Row 136
$ws.Range("H136").Value = 3165.4736
$ws.Range("I136").Value = 1946.4
$ws.Range("K136").Value = 5839.200000000001
$ws.Range("M136").Value = -3289.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 1186
$ws.Range("J18").Value = 1500
$ws.Range("L18").Value = 4500
$ws.Range("N18").Value = -4838
# Row 26
$ws.Range("H26").Value = 2999.5
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 2999.5
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 8998.5
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -9574.5
# Row 34
$ws.Range("H34").Value = 3212.5715
$ws.Range("I34").Value = 1499.3334
$ws.Range("J34").Value = 4497.5
$ws.Range("K34").Value = 4498.0002
$ws.Range("L34").Value = 13492.5
$ws.Range("M34").Value = -4414.0002
$ws.Range("N34").Value = -13660.5
# Row 50
$ws.Range("H50").Value = 578.3333
$ws.Range("I50").Value = 800
$ws.Range("J50").Value = 534
$ws.Range("K50").Value = 2400
$ws.Range("L50").Value = 1602
$ws.Range("M50").Value = -1919
$ws.Range("N50").Value = -2564
# Row 53
$ws.Range("H53").Value = 578.3333
$ws.Range("I53").Value = 800
$ws.Range("J53").Value = 534
$ws.Range("K53").Value = 2400
$ws.Range("L53").Value = 1602
$ws.Range("M53").Value = -1919
$ws.Range("N53").Value = -2564
# Row 57
$ws.Range("H57").Value = 3516.3333
$ws.Range("I57").Value = 3399.5
$ws.Range("K57").Value = 10198.5
$ws.Range("M57").Value = -9639.5
# Row 60
$ws.Range("H60").Value = 4028.923
$ws.Range("I60").Value = 1905.2
$ws.Range("J60").Value = 5356.25
$ws.Range("K60").Value = 5715.6
$ws.Range("L60").Value = 16068.75
$ws.Range("M60").Value = -5464.6
$ws.Range("N60").Value = -16570.75
# Row 61
$ws.Range("H61").Value = 3448.875
$ws.Range("I61").Value = 147
$ws.Range("J61").Value = 6750.75
$ws.Range("K61").Value = 441
$ws.Range("L61").Value = 20252.25
$ws.Range("M61").Value = -226
$ws.Range("N61").Value = -20682.25
# Row 96
$ws.Range("H96").Value = 5500
$ws.Range("J96").Value = 5500
$ws.Range("L96").Value = 16500
$ws.Range("N96").Value = -20618
# Row 100
$ws.Range("H100").Value = 7499
$ws.Range("J100").Value = 7499
$ws.Range("L100").Value = 22497
$ws.Range("N100").Value = -24119
# Row 104
$ws.Range("H104").Value = 6166.6665
$ws.Range("J104").Value = 6166.6665
$ws.Range("L104").Value = 18499.9995
$ws.Range("N104").Value = -23741.9995
# Row 115
$ws.Range("H115").Value = 9000
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 8148.0835
$ws.Range("I132").Value = 2064.6667
$ws.Range("K132").Value = 6194.000100000001
$ws.Range("M132").Value = -3664.000100000001
# Row 136
$ws.Range("H136").Value = 9397.933999999999
$ws.Range("J136").Value = 9397.933999999999
$ws.Range("L136").Value = 28193.802
$ws.Range("N136").Value = -33293.802

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4080.5667
$ws.Range("I40").Value = 3993.4443
$ws.Range("K40").Value = 3993.4443
$ws.Range("M40").Value = -3857.4443
# Row 46
$ws.Range("H46").Value = 1750.2858
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 2150.4
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 2150.4
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -2526.4
# Row 47
$ws.Range("H47").Value = 39500
$ws.Range("J47").Value = 39500
$ws.Range("L47").Value = 39500
$ws.Range("N47").Value = -40480
# Row 52
$ws.Range("H52").Value = 39500
$ws.Range("J52").Value = 39500
$ws.Range("L52").Value = 39500
$ws.Range("N52").Value = -39966
# Row 93
$ws.Range("H93").Value = 2861.4
$ws.Range("I93").Value = 2861.4
$ws.Range("K93").Value = 2861.4
$ws.Range("M93").Value = -1613.4
# Row 122
$ws.Range("H122").Value = 2333.7222
$ws.Range("I122").Value = 2229.1428
$ws.Range("J122").Value = 2699.75
$ws.Range("K122").Value = 6687.428400000001
$ws.Range("L122").Value = 8099.25
$ws.Range("M122").Value = -4237.428400000001
$ws.Range("N122").Value = -12999.25
# Row 132
$ws.Range("H132").Value = 6725.125
$ws.Range("J132").Value = 12777
$ws.Range("L132").Value = 38331
$ws.Range("N132").Value = -43391
